$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.833.71"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "'1.633.27"
$ws.Range("D3").ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'214.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'0.507"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.257"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "'19.90"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "'1.859.52"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "'1.629.73"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "'63.00"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "'25.836.56"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "'193.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").Value = "'9.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").Value = "'6.17"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  -5.14%  "
$ws.Range("D26").Value = "'138.90"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("E27").Value = "  -4.17%  "
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "'15.51"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").Value = "'0.899"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").Value = "'1.120.65"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("E41").Value = "  +1.14%  "
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").Value = "'99.40"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").Value = "'0.0₆0110"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("D46").Value = "'55.34"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("E47").Value = "  -4.67%  "
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "'2.34"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.17%  "
$ws.Range("B51").Value = "Frax"
$ws.Range("C51").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("E51").Value = "  -0.35%  "

Write-Host "Applied 74 cell updates"
